$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hide row 11 (existing data row stays the same, only becomes hidden) ---
$ws.Rows.Item(11).Hidden = $true

# --- Fill rows 12-16 with the same look & feel (borders/alignment) as row 11 ---
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G16").PasteSpecial(-4122)

# Row 12
$ws.Range("A12").Value = 45026
$ws.Range("B12").Value = "Running dan Fixing Script"
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = "Tercapai"
$ws.Range("G12").Value = "1. LinkAja Positive`n2. Tarik Tunai Tanpa Kartu Indomaret Negative`n3. Tarik Tunai Tanpa Kartu ATM BSI Negative`n4. Paket Data XL  Negative `n5. Flow Favourite Transaction`n"
$ws.Rows.Item(12).RowHeight = 120
$ws.Rows.Item(12).Hidden = $true

# Row 13 (note: Keterangan (G) is filled in before Nama (B) so new shared
# strings land in the same order the original authoring session produced)
$ws.Range("A13").Value = 45027
$ws.Range("G13").Value = "1. Tarik Tunai Indomaret Positive`n2. Tarik Tunai ATM BSI Positive`n3. Link Aja Syariah Positive + Negative`n4. LinkAja Negative`n5. BPJS Ketenagakerjaan BPU Negative`n6. Kitabisa Positive + Negative`n7. List Recent Transaction Positive`n8. Universitas Ahmad Dahlan Positive + Negative`n9. Jadiberkah Positive Negative"
$ws.Range("B13").Value = "Running SuperApp Rebrand"
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = "Tercapai"
$ws.Rows.Item(13).RowHeight = 150
$ws.Rows.Item(13).Hidden = $true

# Row 14
$ws.Range("A14").Value = 45028
$ws.Range("B14").Value = "Running SuperApp Rebrand"
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = "Tercapai"
$ws.Range("G14").Value = "1. Pelaporan Transaksi`n2. Pelaporan Melalui Mutasi`n3. Akademik Bayar ID Positive + Negative"
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(14).Hidden = $true

# Row 15
$ws.Range("A15").Value = 45029
$ws.Range("B15").Value = "Running SuperApp Rebrand"
$ws.Range("C15").Value = 13
$ws.Range("D15").Value = 13
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = "Tercapai"
$ws.Range("G15").Value = "1. Rumah Zakat Indonesia Positive + Negative`n2. Pulsa Indosat Positive + Negative`n3. Paket Data Telkomsel Positive + Negative`n4. Paket Data Indosat Positive + Negative`n5. Pulsa XL Positive + Negative`n6. Pulsa Tri Positive"
$ws.Rows.Item(15).RowHeight = 135
$ws.Rows.Item(15).Hidden = $true

# Row 16
$ws.Range("A16").Value = 45030
$ws.Range("B16").Value = "Running SuperApp Rebrand"
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = "Tercapai"
$ws.Range("G16").Value = "1. Pulsa Tri Negative`n2. Paket Data Telkomsel Positive + Negative`n3. Paket Data Tri Positive +Negative `n4. Pulsa Smartfren Positive + Negative`n5. Pulsa Telkomsel Positive + Negative`n6. Dana Positive + Negative`n7. Gopay Positive + Negative`n8. Ovo Positive + Negative"
$ws.Rows.Item(16).RowHeight = 135

# Distinct formatting used only for row 16: G16 left/vcenter aligned with wrap+border,
# C16 centered without border (order matters: these mint two brand-new cell styles,
# and the new-style indices must come out G16 first, then C16, to match the workbook)
$ws.Range("G16").HorizontalAlignment = -4131
$ws.Range("G16").VerticalAlignment = -4108

$ws.Range("C16").Borders.LineStyle = -4142
$ws.Range("C16").HorizontalAlignment = -4108
$ws.Range("C16").VerticalAlignment = -4108

# --- Column G got a bit wider to fit the new notes ---
$ws.Columns.Item(7).ColumnWidth = 38.6

# --- Selection moved from F20 to F19 ---
[void]$ws.Range("F19").Select()
